# Fruta / hortaliza, semanal
# Rotates the weekly price-record data (columns D, M, N, O, P, Q, S, T)
# across rows 2-4: row2 <- old row3, row3 <- old row4, row4 <- old row2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (new values = previous row 3 values)
$ws.Range("D2").Value = 44875
$ws.Range("M2").Value = 50

# Row 3 (new values = previous row 4 values)
$ws.Range("D3").Value = 44855
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 5 kilos"
$ws.Range("S3").Value = 3000
$ws.Range("T3").Value = 5

# Row 4 (new values = previous row 2 values)
$ws.Range("D4").Value = 44874
$ws.Range("M4").Value = 67
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("Q4").Value = "$/bandeja 10 kilos"
$ws.Range("S4").Value = 1600
$ws.Range("T4").Value = 10
